$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.549.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.125.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5249"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4477"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09415"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.182"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.754"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.89%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.974"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.090.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06737"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.362"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.007"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.568.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.329"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.373.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.553"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.163"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1064"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.899"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.302"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02657"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7158"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6959"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.393"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.343"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.666"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("E51").Value = "  +3.16%  "
